$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.451.30'
$ws.Range('E2').Value = '  +4.79%  '
$ws.Range('D3').Value = '2.749.33'
$ws.Range('E3').Value = '  +4.57%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '115.62'
$ws.Range('E5').Value = '  +3.88%  '
$ws.Range('D6').Value = '332.54'
$ws.Range('E6').Value = '  +2.92%  '
$ws.Range('D7').Value = '0.537'
$ws.Range('E7').Value = '  +2.21%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.573'
$ws.Range('E9').Value = '  +5.64%  '
$ws.Range('D10').Value = '41.52'
$ws.Range('E10').Value = '  +4.59%  '
$ws.Range('D11').Value = '0.0856'
$ws.Range('E11').Value = '  +5.65%  '
$ws.Range('D12').Value = '20.17'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('E14').Value = '  +4.93%  '
$ws.Range('D15').Value = '3.179.27'
$ws.Range('E15').Value = '  +4.59%  '
$ws.Range('D16').Value = '2.746.01'
$ws.Range('E16').Value = '  +4.66%  '
$ws.Range('D17').Value = '0.882'
$ws.Range('E17').Value = '  +2.90%  '
$ws.Range('D18').Value = '51.449.58'
$ws.Range('E18').Value = '  +4.86%  '
$ws.Range('D19').Value = '3.21'
$ws.Range('E19').Value = '  +7.94%  '
$ws.Range('D20').Value = '13.41'
$ws.Range('E20').Value = '  +4.03%  '
$ws.Range('E21').Value = '  +2.37%  '
$ws.Range('D22').Value = '0.0₃0973'
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('D23').Value = '278.14'
$ws.Range('E23').Value = '  +3.27%  '
$ws.Range('D24').Value = '69.43'
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('E25').Value = '  +4.51%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '10.19'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').Value = '35.03'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').Value = "'49.80"
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').Value = '5.54'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').Value = '0.0823'
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('D35').Value = '19.13'
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = "'5.00"
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('D38').Value = '2.09'
$ws.Range('E38').Value = '  +2.54%  '
$ws.Range('E39').Value = '  +3.01%  '
$ws.Range('D40').Value = '126.94'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('D41').Value = "'23.00"
$ws.Range('E41').Value = '  +3.86%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0344'
$ws.Range('E42').Value = '  +8.47%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '2.29'
$ws.Range('E43').Value = '  +7.74%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.113'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('D45').Value = '2.44'
$ws.Range('E45').Value = '  +13.57%  '
$ws.Range('D46').Value = '2.090.54'
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('D47').Value = '3.32'
$ws.Range('E47').Value = '  +2.35%  '
$ws.Range('D48').Value = '2.22'
$ws.Range('E49').Value = '  +6.38%  '
$ws.Range('D50').Value = '8.97'
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('D51').Value = '59.81'
